$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.968.62"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.774.69"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'356.18"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'108.88"
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'40.11"
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "'19.34"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("D14").Value = "'7.58"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "3.210.97"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "2.769.92"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "'0.927"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "51.844.06"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'7.36"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "'273.92"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'69.61"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'2.74"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "'26.51"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("D32").Value = "'51.40"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "'33.76"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("D35").Value = "'5.31"
$ws.Range("E35").Value = "  +9.20%  "
$ws.Range("D36").Value = "'0.0835"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'3.19"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "'18.18"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "'122.74"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").Value = "'21.83"
$ws.Range("E45").Value = "  -6.46%  "
$ws.Range("D46").Value = "2.058.71"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'3.24"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'0.923"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E51").Value = "  +0.10%  "
